$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 131, shifting rows 131:153 down to 132:154
$ws.Range("A131:R131").Insert($xlShiftDown)

# Populate the newly inserted row 131 with the new weekly record
$ws.Range("A131").Value = 10
$ws.Range("B131").Value = "Vega Modelo de Temuco"
$ws.Range("C131").Value = "La Araucanía"
$ws.Range("D131").Value = 44474
$ws.Range("E131").Value = 9
$ws.Range("F131").Value = 100112039
$ws.Range("G131").Value = "Ciboulette"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 30
$ws.Range("K131").Value = 5000
$ws.Range("L131").Value = 5000
$ws.Range("M131").Value = 5000
$ws.Range("N131").Value = "$/docena de atados"
$ws.Range("O131").Value = "Región Metropolitana"
$ws.Range("P131").Value = 1667
$ws.Range("Q131").Value = 3
$ws.Range("R131").Value = "Hortaliza"

Write-Host "done"
